$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark all Tc (test cases) to Yes: update C2 and C3 from "N" to "Y"
$ws.Range("C2").Value = "Y"
$ws.Range("C3").Value = "Y"

# Update the active selection to C2, matching the recorded cursor position
$ws.Range("C2").Select()
